# Insert a new weekly record at row 619 ("Fruta / hortaliza, semanal").
# This shifts the existing rows 619:650 down to 620:651, then the new
# row 619 is populated with the same record as the (now shifted) row 620,
# except for an updated date (the new weekly observation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 619; everything below moves down by one row.
$ws.Rows.Item(619).Insert()

# Duplicate the row that is now at 620 (the old row 619) into the new
# blank row 619 so all of its formatting/values are restored.
$ws.Range("A620:R620").Copy()
$ws.Range("A619:R619").PasteSpecial()

# The new row represents a newer weekly observation (2023-08-09 instead
# of 2023-07-12), so only the date changes relative to row 620.
$ws.Cells.Item(619, 4).Value = 45147

$excel.CutCopyMode = 0
